$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-09-13 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-14 Thursday", 2)

$d.Content.Find.Execute("31×59=1829", $true, $false, $false, $false, $false, $true, 1, $false, "40×76=3040", 2)
$d.Content.Find.Execute("53×48=2544", $true, $false, $false, $false, $false, $true, 1, $false, "65×85=5525", 2)
$d.Content.Find.Execute("92×52=4784", $true, $false, $false, $false, $false, $true, 1, $false, "47×92=4324", 2)
$d.Content.Find.Execute("65×21=1365", $true, $false, $false, $false, $false, $true, 1, $false, "46×74=3404", 2)
$d.Content.Find.Execute("33×69=2277", $true, $false, $false, $false, $false, $true, 1, $false, "78×60=4680", 2)

$d.Content.Find.Execute("70×92=6440", $true, $false, $false, $false, $false, $true, 1, $false, "52×74=3848", 2)
$d.Content.Find.Execute("33×36=1188", $true, $false, $false, $false, $false, $true, 1, $false, "40×79=3160", 2)
$d.Content.Find.Execute("79×46=3634", $true, $false, $false, $false, $false, $true, 1, $false, "77×88=6776", 2)
$d.Content.Find.Execute("37×76=2812", $true, $false, $false, $false, $false, $true, 1, $false, "33×45=1485", 2)
$d.Content.Find.Execute("50×81=4050", $true, $false, $false, $false, $false, $true, 1, $false, "36×57=2052", 2)

$d.Content.Find.Execute("35×81=2835", $true, $false, $false, $false, $false, $true, 1, $false, "61×66=4026", 2)
$d.Content.Find.Execute("28×94=2632", $true, $false, $false, $false, $false, $true, 1, $false, "27×77=2079", 2)
$d.Content.Find.Execute("22×64=1408", $true, $false, $false, $false, $false, $true, 1, $false, "70×36=2520", 2)
$d.Content.Find.Execute("31×60=1860", $true, $false, $false, $false, $false, $true, 1, $false, "22×50=1100", 2)
$d.Content.Find.Execute("50×77=3850", $true, $false, $false, $false, $false, $true, 1, $false, "96×90=8640", 2)

$d.Content.Find.Execute("95×57=5415", $true, $false, $false, $false, $false, $true, 1, $false, "81×53=4293", 2)
$d.Content.Find.Execute("92×33=3036", $true, $false, $false, $false, $false, $true, 1, $false, "56×64=3584", 2)
$d.Content.Find.Execute("12×16=192", $true, $false, $false, $false, $false, $true, 1, $false, "85×28=2380", 2)
$d.Content.Find.Execute("80×12=960", $true, $false, $false, $false, $false, $true, 1, $false, "40×11=440", 2)
$d.Content.Find.Execute("12×26=312", $true, $false, $false, $false, $false, $true, 1, $false, "95×32=3040", 2)

$d.Content.Find.Execute("79×20=1580", $true, $false, $false, $false, $false, $true, 1, $false, "12×65=780", 2)
$d.Content.Find.Execute("48×19=912", $true, $false, $false, $false, $false, $true, 1, $false, "40×31=1240", 2)
$d.Content.Find.Execute("78×39=3042", $true, $false, $false, $false, $false, $true, 1, $false, "61×42=2562", 2)
$d.Content.Find.Execute("85×98=8330", $true, $false, $false, $false, $false, $true, 1, $false, "86×41=3526", 2)
$d.Content.Find.Execute("42×65=2730", $true, $false, $false, $false, $false, $true, 1, $false, "54×99=5346", 2)
